# Weekly update: prepend a new week of "Murcott" mandarina price rows
# for Comercializadora del Agro de Limarí (Coquimbo), shifting the
# existing historical rows down by 3 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the existing data block (old row 428).
$ws.Range("428:430").Insert()

# Row 428: Murcott / Especial
$ws.Cells.Item(428, 1).Value = 2
$ws.Cells.Item(428, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(428, 3).Value = "Coquimbo"
$ws.Cells.Item(428, 4).Value = 44826
$ws.Cells.Item(428, 5).Value = 4
$ws.Cells.Item(428, 6).Value = "Fruta"
$ws.Cells.Item(428, 7).Value = 100102
$ws.Cells.Item(428, 8).Value = "Cítricos"
$ws.Cells.Item(428, 9).Value = 100102004
$ws.Cells.Item(428, 10).Value = "Mandarina"
$ws.Cells.Item(428, 11).Value = "Murcott"
$ws.Cells.Item(428, 12).Value = "Especial"
$ws.Cells.Item(428, 13).Value = 520
$ws.Cells.Item(428, 14).Value = 5500
$ws.Cells.Item(428, 15).Value = 6000
$ws.Cells.Item(428, 16).Value = 5750
$ws.Cells.Item(428, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(428, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(428, 19).Value = 575
$ws.Cells.Item(428, 20).Value = 10

# Row 429: Murcott / Primera
$ws.Cells.Item(429, 1).Value = 2
$ws.Cells.Item(429, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(429, 3).Value = "Coquimbo"
$ws.Cells.Item(429, 4).Value = 44826
$ws.Cells.Item(429, 5).Value = 4
$ws.Cells.Item(429, 6).Value = "Fruta"
$ws.Cells.Item(429, 7).Value = 100102
$ws.Cells.Item(429, 8).Value = "Cítricos"
$ws.Cells.Item(429, 9).Value = 100102004
$ws.Cells.Item(429, 10).Value = "Mandarina"
$ws.Cells.Item(429, 11).Value = "Murcott"
$ws.Cells.Item(429, 12).Value = "Primera"
$ws.Cells.Item(429, 13).Value = 360
$ws.Cells.Item(429, 14).Value = 4500
$ws.Cells.Item(429, 15).Value = 5000
$ws.Cells.Item(429, 16).Value = 4750
$ws.Cells.Item(429, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(429, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(429, 19).Value = 475
$ws.Cells.Item(429, 20).Value = 10

# Row 430: Murcott / Segunda
$ws.Cells.Item(430, 1).Value = 2
$ws.Cells.Item(430, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(430, 3).Value = "Coquimbo"
$ws.Cells.Item(430, 4).Value = 44826
$ws.Cells.Item(430, 5).Value = 4
$ws.Cells.Item(430, 6).Value = "Fruta"
$ws.Cells.Item(430, 7).Value = 100102
$ws.Cells.Item(430, 8).Value = "Cítricos"
$ws.Cells.Item(430, 9).Value = 100102004
$ws.Cells.Item(430, 10).Value = "Mandarina"
$ws.Cells.Item(430, 11).Value = "Murcott"
$ws.Cells.Item(430, 12).Value = "Segunda"
$ws.Cells.Item(430, 13).Value = 300
$ws.Cells.Item(430, 14).Value = 3500
$ws.Cells.Item(430, 15).Value = 4000
$ws.Cells.Item(430, 16).Value = 3750
$ws.Cells.Item(430, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(430, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(430, 19).Value = 375
$ws.Cells.Item(430, 20).Value = 10
